$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Switch workbook calculation mode to manual (calcPr calcMode="manual")
$excel.Calculation = -4135

# Add new row 11 with the week's update (date, name, item, goal, status)
# Copy formatting from analogous existing cells, then set the values afterwards
# so the shared-string table order matches: C11, E11, D11.

# A11: date cell, same format as A9/A10 (style used for date column)
$ws.Range("A9").Copy()
$ws.Range("A11").PasteSpecial(-4122)
$ws.Range("A11").Value = 43824

# B11: name cell, vertically centered like B7
$ws.Range("B7").Copy()
$ws.Range("B11").PasteSpecial(-4122)
$ws.Range("B11").Value = "楊雅婷"

# C11: work item
$ws.Range("C11").Value = "製作期末影片"

# E11: completion status
$ws.Range("E4").Copy()
$ws.Range("E11").PasteSpecial(-4122)
$ws.Range("E11").Value = "預計今日晚上完成"

# D11: completion goal
$ws.Range("D4").Copy()
$ws.Range("D11").PasteSpecial(-4122)
$ws.Range("D11").Value = "完成期末影片"

# Update the active selection to D14, matching the saved cursor position
[void]$ws.Range("D14").Select()
